# Hortaliza, Vega Monumental Concepción - Cilantro
# Weekly refresh: insert the new week's two rows (Primera / Segunda) at the
# top of the data block (row 170) and let everything below shift down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 170, pushing the existing rows 170..271 down
# to 172..273 (dimension grows from A1:R271 to A1:R273 automatically).
$ws.Rows.Item(170).Insert()
$ws.Rows.Item(170).Insert()

# New row 170: Cilantro, Primera, new week (2023-02-03)
$ws.Cells.Item(170, 1).Value = 11
$ws.Cells.Item(170, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(170, 3).Value = "Bíobío"
$ws.Cells.Item(170, 4).Value = 44960
$ws.Cells.Item(170, 5).Value = 8
$ws.Cells.Item(170, 6).Value = 100112040
$ws.Cells.Item(170, 7).Value = "Cilantro"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 100
$ws.Cells.Item(170, 11).Value = 700
$ws.Cells.Item(170, 12).Value = 800
$ws.Cells.Item(170, 13).Value = 750
$ws.Cells.Item(170, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(170, 15).Value = "Región de Ñuble"
$ws.Cells.Item(170, 16).Value = 750
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# New row 171: Cilantro, Segunda, same new week
$ws.Cells.Item(171, 1).Value = 11
$ws.Cells.Item(171, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(171, 3).Value = "Bíobío"
$ws.Cells.Item(171, 4).Value = 44960
$ws.Cells.Item(171, 5).Value = 8
$ws.Cells.Item(171, 6).Value = 100112040
$ws.Cells.Item(171, 7).Value = "Cilantro"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Segunda"
$ws.Cells.Item(171, 10).Value = 50
$ws.Cells.Item(171, 11).Value = 600
$ws.Cells.Item(171, 12).Value = 600
$ws.Cells.Item(171, 13).Value = 600
$ws.Cells.Item(171, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(171, 15).Value = "Región de Ñuble"
$ws.Cells.Item(171, 16).Value = 600
$ws.Cells.Item(171, 17).Value = 1
$ws.Cells.Item(171, 18).Value = "Hortaliza"
